$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: the registration number was previously prefixed with
# "Último registro encontrado: " and status was "Pendente".
# After adjusting the XPath for ANVISA's new site structure, the
# scraper now returns just the bare registration number and the
# status becomes "OK" (matching the already-correct row 5 pattern).
foreach ($row in 2..4) {
    # Force text format so the numeric-looking registration code
    # is stored as a string, not a number (consistent with row 5).
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = "104910019"
    $ws.Range("F$row").Value = "OK"
}
